# Feria Lagunitas de Puerto Montt - Pepino ensalada
# Weekly data refresh: insert a new latest-week observation at row 220,
# pushing the existing historical rows (220-339) down by one (221-340).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 220:339 down to 221:340, leaving row 220 blank for the new entry.
$ws.Rows("220:220").Insert()

# Populate the newly inserted row 220 with the new week's record.
$ws.Cells.Item(220, 1).Value = 4
$ws.Cells.Item(220, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(220, 3).Value = "Los Lagos"
$ws.Cells.Item(220, 4).Value = 44873
$ws.Cells.Item(220, 5).Value = 10
$ws.Cells.Item(220, 6).Value = 100112043
$ws.Cells.Item(220, 7).Value = "Pepino ensalada"
$ws.Cells.Item(220, 8).Value = "Sin especificar"
$ws.Cells.Item(220, 9).Value = "Primera"
$ws.Cells.Item(220, 10).Value = 400
$ws.Cells.Item(220, 11).Value = 22000
$ws.Cells.Item(220, 12).Value = 23000
$ws.Cells.Item(220, 13).Value = 22500
$ws.Cells.Item(220, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(220, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(220, 16).Value = 375
$ws.Cells.Item(220, 17).Value = 60
$ws.Cells.Item(220, 18).Value = "Hortaliza"
